# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a handful of rows on the single
# worksheet in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (new DAMSLTag, new DialogAct)
$updates = @{
    23 = @("aa", "Agree/Accept")
    28 = @("sd", "Statement-non-opinion")
    37 = @("sv", "Statement-opinion")
    39 = @("aa", "Agree/Accept")
    51 = @("b",  "Acknowledge (Backchannel)")
    68 = @("aa", "Agree/Accept")
    80 = @("sv", "Statement-opinion")
    82 = @("aa", "Agree/Accept")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}
